$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the existing column C (T.C (Azure)) to hold
# the new "PRINCIPAL" field, shifting T.C (Azure), T.C (Desc.) and Error
# one column to the right.
$ws.Range("C:C").Insert()

# New header + value for the inserted PRINCIPAL column.
$ws.Range("C1").Value = "PRINCIPAL"
$ws.Range("C2").Value = 10
